$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the two new daily rows (2025-09-01 / serial 45901) for each station.
$ws.Range("A64").Value = 45901
$ws.Range("B64").Value = "四方坪站"
$ws.Range("C64").Value = 11665.57
$ws.Range("D64").Value = 9409.17
$ws.Range("E64").Value = 4048.22
$ws.Range("F64").Value = 466

$ws.Range("A65").Value = 45901
$ws.Range("B65").Value = "高岭站"
$ws.Range("C65").Value = 4787.61
$ws.Range("D65").Value = 3921.98
$ws.Range("E65").Value = 1197.92
$ws.Range("F65").Value = 164

# Mirror the author's updated selection (scrolled down two rows further).
$ws.Range("H62").Select() | Out-Null
